# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (new rows 279-281) above the existing
# data block, pushing the former rows 279-337 down to rows 282-340.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the top of the data block (rows 279:281).
# Everything that used to live at row 279 onward shifts down by 3 rows,
# carrying its formatting (incl. the date column's number format) with it.
$ws.Rows("279:281").Insert()

# --- New row 279 -----------------------------------------------------
$ws.Cells.Item(279,1).Value = 10
$ws.Cells.Item(279,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(279,3).Value = "La Araucanía"
$ws.Cells.Item(279,4).Value = 45218
$ws.Cells.Item(279,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(279,5).Value = 9
$ws.Cells.Item(279,6).Value = "Fruta"
$ws.Cells.Item(279,7).Value = 100101
$ws.Cells.Item(279,8).Value = "Berries"
$ws.Cells.Item(279,9).Value = 100112025
$ws.Cells.Item(279,10).Value = "Frutilla"
$ws.Cells.Item(279,11).Value = "Sin especificar"
$ws.Cells.Item(279,12).Value = "Primera"
$ws.Cells.Item(279,13).Value = 1800
$ws.Cells.Item(279,14).Value = 10000
$ws.Cells.Item(279,15).Value = 11000
$ws.Cells.Item(279,16).Value = 10556
$ws.Cells.Item(279,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(279,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(279,19).Value = 1508
$ws.Cells.Item(279,20).Value = 7

# --- New row 280 -----------------------------------------------------
$ws.Cells.Item(280,1).Value = 10
$ws.Cells.Item(280,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(280,3).Value = "La Araucanía"
$ws.Cells.Item(280,4).Value = 45218
$ws.Cells.Item(280,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(280,5).Value = 9
$ws.Cells.Item(280,6).Value = "Fruta"
$ws.Cells.Item(280,7).Value = 100101
$ws.Cells.Item(280,8).Value = "Berries"
$ws.Cells.Item(280,9).Value = 100112025
$ws.Cells.Item(280,10).Value = "Frutilla"
$ws.Cells.Item(280,11).Value = "Sin especificar"
$ws.Cells.Item(280,12).Value = "Segunda"
$ws.Cells.Item(280,13).Value = 330
$ws.Cells.Item(280,14).Value = 7000
$ws.Cells.Item(280,15).Value = 8000
$ws.Cells.Item(280,16).Value = 7545
$ws.Cells.Item(280,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(280,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(280,19).Value = 1078
$ws.Cells.Item(280,20).Value = 7

# --- New row 281 -----------------------------------------------------
$ws.Cells.Item(281,1).Value = 10
$ws.Cells.Item(281,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(281,3).Value = "La Araucanía"
$ws.Cells.Item(281,4).Value = 45218
$ws.Cells.Item(281,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(281,5).Value = 9
$ws.Cells.Item(281,6).Value = "Fruta"
$ws.Cells.Item(281,7).Value = 100101
$ws.Cells.Item(281,8).Value = "Berries"
$ws.Cells.Item(281,9).Value = 100112025
$ws.Cells.Item(281,10).Value = "Frutilla"
$ws.Cells.Item(281,11).Value = "Sin especificar"
$ws.Cells.Item(281,12).Value = "Tercera"
$ws.Cells.Item(281,13).Value = 100
$ws.Cells.Item(281,14).Value = 6000
$ws.Cells.Item(281,15).Value = 6000
$ws.Cells.Item(281,16).Value = 6000
$ws.Cells.Item(281,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(281,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(281,19).Value = 857
$ws.Cells.Item(281,20).Value = 7
